$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "27.198.78"
Set-TextValue $ws.Cells.Item(2, 5) "  -1.92%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "1.817.62"
Set-TextValue $ws.Cells.Item(3, 5) "  -2.48%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 5) "  -1.40%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "313.67"
Set-TextValue $ws.Cells.Item(5, 5) "  -2.28%  "

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "1.004"
Set-TextValue $ws.Cells.Item(6, 5) "  -1.53%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.4259"
Set-TextValue $ws.Cells.Item(7, 5) "  -2.51%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.3661"
Set-TextValue $ws.Cells.Item(8, 5) "  -3.62%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "45.79"
Set-TextValue $ws.Cells.Item(9, 5) "  -2.07%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "0.07193"
Set-TextValue $ws.Cells.Item(10, 5) "  -3.33%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.8594"
Set-TextValue $ws.Cells.Item(11, 5) "  -2.85%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "20.94"
Set-TextValue $ws.Cells.Item(12, 5) "  -3.20%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "1.825.58"
Set-TextValue $ws.Cells.Item(13, 5) "  -2.14%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "6.642"
Set-TextValue $ws.Cells.Item(14, 5) "  -1.69%  "

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "0.07080"
Set-TextValue $ws.Cells.Item(15, 5) "  -0.39%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "5.296"
Set-TextValue $ws.Cells.Item(16, 5) "  -3.65%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "87.75"
Set-TextValue $ws.Cells.Item(17, 5) "  +1.31%  "

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "1.006"
Set-TextValue $ws.Cells.Item(18, 5) "  -1.87%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "0.000008852"
Set-TextValue $ws.Cells.Item(19, 5) "  -2.68%  "

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "1.004"
Set-TextValue $ws.Cells.Item(20, 5) "  -1.60%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "15.00"
Set-TextValue $ws.Cells.Item(21, 5) "  -3.09%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "27.227.76"
Set-TextValue $ws.Cells.Item(22, 5) "  -1.79%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "5.110"
Set-TextValue $ws.Cells.Item(23, 5) "  -3.45%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "10.85"
Set-TextValue $ws.Cells.Item(24, 5) "  -2.78%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "2.094.85"
Set-TextValue $ws.Cells.Item(25, 5) "  +0.22%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "2.001"
Set-TextValue $ws.Cells.Item(26, 5) "  -1.99%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "152.84"
Set-TextValue $ws.Cells.Item(27, 5) "  -3.02%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "18.25"
Set-TextValue $ws.Cells.Item(28, 5) "  -2.68%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "2.109"
Set-TextValue $ws.Cells.Item(29, 5) "  +5.48%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "5.220"
Set-TextValue $ws.Cells.Item(30, 5) "  -2.72%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "116.08"
Set-TextValue $ws.Cells.Item(31, 5) "  -3.99%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "0.08873"
Set-TextValue $ws.Cells.Item(32, 5) "  -2.00%  "

# Row 33
Set-TextValue $ws.Cells.Item(33, 2) "ARBITRUM"
Set-TextValue $ws.Cells.Item(33, 3) "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Cells.Item(33, 4) "1.191"
Set-TextValue $ws.Cells.Item(33, 5) "  -2.63%  "

# Row 34
Set-TextValue $ws.Cells.Item(34, 2) "ImmutableX"
Set-TextValue $ws.Cells.Item(34, 3) "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Cells.Item(34, 4) "0.7581"
Set-TextValue $ws.Cells.Item(34, 5) "  -1.30%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "4.454"
Set-TextValue $ws.Cells.Item(35, 5) "  -2.45%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "2.825"
Set-TextValue $ws.Cells.Item(36, 5) "  -6.98%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "1.004"
Set-TextValue $ws.Cells.Item(37, 5) "  -1.72%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "1.113"
Set-TextValue $ws.Cells.Item(38, 5) "  -2.49%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "0.01952"
Set-TextValue $ws.Cells.Item(39, 5) "  -1.49%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.05239"
Set-TextValue $ws.Cells.Item(40, 5) "  -1.08%  "

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "2.897"
Set-TextValue $ws.Cells.Item(41, 5) "  +0.61%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "7.035"
Set-TextValue $ws.Cells.Item(42, 5) "  +1.26%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "0.1673"

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.4997"
Set-TextValue $ws.Cells.Item(44, 5) "  -3.85%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "8.592"
Set-TextValue $ws.Cells.Item(45, 5) "  -1.21%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "10.50"
Set-TextValue $ws.Cells.Item(46, 5) "  -1.99%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "106.16"
Set-TextValue $ws.Cells.Item(47, 5) "  -3.53%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "0.4686"
Set-TextValue $ws.Cells.Item(48, 5) "  -0.75%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 2) "Cronos"
Set-TextValue $ws.Cells.Item(49, 3) "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(49, 4) "0.06400"
Set-TextValue $ws.Cells.Item(49, 5) "  -1.60%  "

# Row 50
Set-TextValue $ws.Cells.Item(50, 2) "PaxDollar"
Set-TextValue $ws.Cells.Item(50, 3) "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Cells.Item(50, 4) "1.003"
Set-TextValue $ws.Cells.Item(50, 5) "  -1.80%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "1.653"
Set-TextValue $ws.Cells.Item(51, 5) "  -3.64%  "
